$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Helper: find the 1-based Paragraphs() index of the first paragraph whose
# text matches a -like pattern. Re-scanned fresh every call so it is never
# stale after an earlier InsertXML shifts indices around.
function Find-ParaIndex($pattern) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text -like $pattern) {
            return $i
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1) "Newt window Pipeline releaseCtx gives error" paragraph becomes two new
#    list items: "Fallout4 BloodySawBat In loadart folder locks the world up"
#    and "Noramls for verts obviously way!! Wrong in fallout4".
#    Do this first (bottom-most of the trio) so it does not disturb the
#    paragraph indices of the two list items above it.
# ---------------------------------------------------------------------------
$xmlNewt = @"
<w:p $wNs><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Fallout4 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>BloodySawBat</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> In </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>loadart</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> folder locks the world up</w:t></w:r></w:p><w:p $wNs><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Noramls</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> for verts obviously way!! Wrong in fallout4</w:t></w:r></w:p>
"@

$idx = Find-ParaIndex("Newt window Pipeline releaseCtx gives error*")
$d.Paragraphs($idx).Range.InsertXML($xmlNewt)
Write-Output "Step1 done"

# ---------------------------------------------------------------------------
# 2) "Free all the declared buffers ... system (or not?)" paragraph (3 runs)
#    becomes "Newt window Pipeline releaseCtx gives error" (spell-checked
#    run around releaseCtx).
# ---------------------------------------------------------------------------
$xmlFree = @"
<w:p $wNs><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Newt window Pipeline </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>releaseCtx</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> gives error</w:t></w:r></w:p>
"@

$idx = Find-ParaIndex("Free all the declared buffers when geo is freed*")
$d.Paragraphs($idx).Range.InsertXML($xmlFree)
Write-Output "Step2 done"

# ---------------------------------------------------------------------------
# 3) "Now using GL2ES2 profile depth buffer returns 0 (but appears to work)"
#    becomes "Free all the declared buffers when geo is freed in a manner
#    that is much better than the current system (or not?)" (single run).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Now using GL2ES2 profile depth buffer returns 0 (but appears to work)", $true, $false, $false, $false, $false, $true, 1, $false, "Free all the declared buffers when geo is freed in a manner that is much better than the current system (or not?)", 2) | Out-Null
Write-Output "Step3 done"

# ---------------------------------------------------------------------------
# 4) "... stereo) and it'll do it's best" paragraph (Consolas run): split the
#    trailing run into three runs with a spell-checked "it's". InsertXML only
#    behaves correctly when given the *whole* paragraph (incl. <w:pPr/>), so
#    rebuild the full paragraph rather than touching a sub-range.
# ---------------------------------------------------------------------------
$rPr = '<w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
$apos = [char]0x2019

$pPrStereo = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/>' + $rPr + '</w:pPr>'

$runsBefore = `
  '<w:r>' + $rPr + '<w:t xml:space="preserve">All the get graphics config gear chucked, </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r>' + $rPr + '<w:t>graphicconfiguration</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r>' + $rPr + '<w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r>' + $rPr + '<w:t>graphicconfigtemple</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r>' + $rPr + '<w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r>' + $rPr + '<w:t>etc</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r>' + $rPr + '<w:t xml:space="preserve"> all gone. Just get the Canvas3D what you what to have(</w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r>' + $rPr + '<w:t>color</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r>' + $rPr + '<w:t xml:space="preserve">, depth, stencil, </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r>' + $rPr + '<w:t>multisample</w:t></w:r><w:proofErr w:type="spellEnd"/>'

$runsAfter = `
  '<w:r>' + $rPr + '<w:t xml:space="preserve">, stereo) and it' + $apos + 'll do </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/><w:r>' + $rPr + '<w:t>it' + $apos + 's</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
  '<w:r>' + $rPr + '<w:t xml:space="preserve"> best</w:t></w:r>'

$xmlStereo = '<w:p ' + $wNs + '>' + $pPrStereo + $runsBefore + $runsAfter + '</w:p>'

# ---------------------------------------------------------------------------
# 5) New list item right after it: "View. VIRTUAL_WORLD mode needs to be
#    removed, god know what it does now" (same Consolas list, numId 1).
# ---------------------------------------------------------------------------
$pPrView = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/>' + $rPr + '</w:pPr>'
$runsView = `
  '<w:r>' + $rPr + '<w:t>View.</w:t></w:r>' + `
  '<w:r>' + $rPr + '<w:t xml:space="preserve"> </w:t></w:r>' + `
  '<w:r>' + $rPr + '<w:t>VIRTUAL_WORLD</w:t></w:r>' + `
  '<w:r>' + $rPr + '<w:t xml:space="preserve"> mode needs to be removed, god know what it does now</w:t></w:r>'
$xmlView = '<w:p ' + $wNs + '>' + $pPrView + $runsView + '</w:p>'

$idx = Find-ParaIndex("*stereo) and it*ll do it*s best*")
$d.Paragraphs($idx).Range.InsertXML($xmlStereo + $xmlView)
Write-Output "Step4+5 done"
